$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected (with a password) - unprotect temporarily so the
# cell values below can be updated, then restore protection afterwards.
$ws.Unprotect("lido")

# Update confidential disclaimer text (date 2021-05-13 -> 2021-05-14)
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-39
$ws.Range("D2").Value = 0.05758430141939996
$ws.Range("E2").Value = 0.01984476274305846
$ws.Range("D3").Value = 0.0520858504606939
$ws.Range("E3").Value = 0.021067357939349
$ws.Range("D4").Value = 0.3128730148337406
$ws.Range("E4").Value = 0.004075761208343298
$ws.Range("D5").Value = 0.03387809193432291
$ws.Range("E5").Value = 0.01943083439033111
$ws.Range("D6").Value = 0.03116856897646041
$ws.Range("E6").Value = 0.01430231518727099
$ws.Range("D7").Value = 0.03115120916949289
$ws.Range("E7").Value = 0.01554179566563452
$ws.Range("D8").Value = 0.02914047202172419
$ws.Range("E8").Value = 0.001529771710990779
$ws.Range("D9").Value = 0.0237019231129863
$ws.Range("E9").Value = 0.0092592592592593
$ws.Range("D10").Value = 0.02388623711288835
$ws.Range("E10").Value = 0.02213508954527521
$ws.Range("D11").Value = 0.02289801402736698
$ws.Range("E11").Value = 0.03498656882657425
$ws.Range("D12").Value = 0.02331443507598291
$ws.Range("E12").Value = 0.01242829827915837
$ws.Range("D13").Value = 0.02102186896324768
$ws.Range("E13").Value = -0.02601771896377714
$ws.Range("D14").Value = 0.02213696866265512
$ws.Range("E14").Value = 0.002105721754283696
$ws.Range("D15").Value = 0.02092306808655599
$ws.Range("E15").Value = -0.005500583860856723
$ws.Range("D16").Value = 0.02189050226743725
$ws.Range("E16").Value = 0.003035049931466638
$ws.Range("D17").Value = 0.01928974600137742
$ws.Range("E17").Value = 0.01080495528026226
$ws.Range("D18").Value = 0.01389041739726986
$ws.Range("E18").Value = 0.02481022032956859
$ws.Range("D19").Value = 0.01721964260015188
$ws.Range("E19").Value = 0.02247778358599062
$ws.Range("D20").Value = 0.01575509641604642
$ws.Range("E20").Value = -0.002040469307940773
$ws.Range("D21").Value = 0.01652182122377852
$ws.Range("E21").Value = 0.02478920741989876
$ws.Range("D22").Value = 0.01225238030279147
$ws.Range("E22").Value = 0.03157305532718779
$ws.Range("D23").Value = 0.01518725927332491
$ws.Range("E23").Value = 0.004035956705191746
$ws.Range("D24").Value = 0.01484188484211309
$ws.Range("E24").Value = 0.0009313877677741278
$ws.Range("D25").Value = 0.01406197944020194
$ws.Range("E25").Value = 0.007811011621261299
$ws.Range("D26").Value = 0.01384444605659661
$ws.Range("E26").Value = 0.007322264793529154
$ws.Range("D27").Value = 0.01288579893850138
$ws.Range("E27").Value = 0.01736396976274213
$ws.Range("D28").Value = 0.01371810523922189
$ws.Range("E28").Value = 0.02615298087739015
$ws.Range("D29").Value = 0.01461009927130603
$ws.Range("E29").Value = -0.001995012468827939
$ws.Range("D30").Value = 0.01339748460683412
$ws.Range("E30").Value = 0.009598157153826392
$ws.Range("D31").Value = 0.01254792565844837
$ws.Range("E31").Value = 0.001503040240486531
$ws.Range("D32").Value = 0.01374425210897543
$ws.Range("E32").Value = -0.001457975986277837
$ws.Range("D33").Value = 0.01264854824327862
$ws.Range("E33").Value = -0.005083238022620229
$ws.Range("D34").Value = 0.00585743462130599
$ws.Range("E34").Value = 0.04227877279961945
$ws.Range("D35").Value = 0.005215014604205509
$ws.Range("E35").Value = 0.01378786010767263
$ws.Range("D36").Value = 0.005160791997257578
$ws.Range("E36").Value = 0.02279900332225915
$ws.Range("D37").Value = 0.005081065476369712
$ws.Range("E37").Value = 0.02615151003880545
$ws.Range("D38").Value = 0.004614279555687527
$ws.Range("E38").Value = 0.02559219693450987
$ws.Range("D39").Value = 0.9999999999999996
$ws.Range("E39").Value = 0.009733493743664612

# Restore the sheet protection that was lifted above.
$ws.Protect("lido", $true, $true, $true, $false, $true, $false, $false, $true, $true, $true, $true, $true, $true, $true, $true)
